$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (bold/border/centered) into H1, then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data column values.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
